$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 392
$ws.Range("J18").Value = 490
$ws.Range("L18").Value = 490
$ws.Range("N18").Value = -1058
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40924
$ws.Range("H112").Value = 19232354
$ws.Range("J112").Value = 1635.0613
$ws.Range("L112").Value = 4905.1839
$ws.Range("N112").Value = -7121.1839
$ws.Range("H125").Value = 3600
$ws.Range("I125").Value = 2666.6667
$ws.Range("J125").Value = 4160
$ws.Range("K125").Value = 24000.0003
$ws.Range("L125").Value = 37440
$ws.Range("M125").Value = -21540.0003
$ws.Range("N125").Value = -42360
$ws.Range("H137").Value = 3381.4707
$ws.Range("I137").Value = 2502.2432
$ws.Range("J137").Value = 5705.143
$ws.Range("K137").Value = 7506.7296
$ws.Range("L137").Value = 17115.429
$ws.Range("M137").Value = -4956.7296
$ws.Range("N137").Value = -22215.429
$ws.Range("H138").Value = 4725.65
$ws.Range("I138").Value = 772.1429000000001
$ws.Range("J138").Value = 5369.244
$ws.Range("K138").Value = 2316.4287
$ws.Range("L138").Value = 16107.732
$ws.Range("M138").Value = 2823.5713
$ws.Range("N138").Value = -26387.732

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1487.23
$ws.Range("I32").Value = 1272.5946
$ws.Range("J32").Value = 2098.1155
$ws.Range("K32").Value = 1272.5946
$ws.Range("L32").Value = 2098.1155
$ws.Range("M32").Value = -985.5945999999999
$ws.Range("N32").Value = -2672.1155
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = 0
$ws.Range("H74").Value = 3057.162
$ws.Range("I74").Value = 3062.9
$ws.Range("K74").Value = 3062.9
$ws.Range("M74").Value = -2188.9
$ws.Range("H77").Value = 3057.162
$ws.Range("I77").Value = 3062.9
$ws.Range("K77").Value = 15314.5
$ws.Range("M77").Value = -10946.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1885.8422
$ws.Range("I134").Value = 1205.5471
$ws.Range("J134").Value = 3453.4783
$ws.Range("K134").Value = 3616.6413
$ws.Range("L134").Value = 10360.4349
$ws.Range("M134").Value = -1081.6413
$ws.Range("N134").Value = -15430.4349

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7814863
$ws.Range("I31").Value = 1334.4043
$ws.Range("K31").Value = 1334.4043
$ws.Range("M31").Value = -1039.4043
$ws.Range("H34").Value = 7814863
$ws.Range("I34").Value = 1334.4043
$ws.Range("K34").Value = 1334.4043
$ws.Range("M34").Value = -1132.4043
$ws.Range("H41").Value = 27395.309
$ws.Range("I41").Value = 8279.5
$ws.Range("K41").Value = 8279.5
$ws.Range("M41").Value = -7851.5
$ws.Range("H122").Value = 1965.8096
$ws.Range("I122").Value = 1076.4445
$ws.Range("J122").Value = 2632.8333
$ws.Range("K122").Value = 3229.3335
$ws.Range("L122").Value = 7898.499899999999
$ws.Range("M122").Value = -779.3335000000002
$ws.Range("N122").Value = -12798.4999
$ws.Range("H132").Value = 3083.4482
$ws.Range("I132").Value = 2889.9473
$ws.Range("J132").Value = 3451.1
$ws.Range("K132").Value = 8669.841899999999
$ws.Range("L132").Value = 10353.3
$ws.Range("M132").Value = -6139.841899999999
$ws.Range("N132").Value = -15413.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6124
$ws.Range("I56").Value = 6124
$ws.Range("K56").Value = 6124
$ws.Range("M56").Value = -5594
$ws.Range("H113").Value = 644.9149
$ws.Range("I113").Value = 533.91174
$ws.Range("J113").Value = 935.2308
$ws.Range("K113").Value = 1601.73522
$ws.Range("L113").Value = 2805.6924
$ws.Range("M113").Value = 568.26478
$ws.Range("N113").Value = -7145.6924

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6501.079
$ws.Range("I70").Value = 5945.615
$ws.Range("J70").Value = 7704.5835
$ws.Range("K70").Value = 5945.615
$ws.Range("L70").Value = 7704.5835
$ws.Range("M70").Value = -5675.615
$ws.Range("N70").Value = -8244.583500000001
$ws.Range("H73").Value = 6501.079
$ws.Range("I73").Value = 5945.615
$ws.Range("J73").Value = 7704.5835
$ws.Range("K73").Value = 5945.615
$ws.Range("L73").Value = 7704.5835
$ws.Range("M73").Value = -5009.615
$ws.Range("N73").Value = -9576.583500000001
$ws.Range("H102").Value = 1990.4117
$ws.Range("I102").Value = 1723.3529
$ws.Range("J102").Value = 2257.4707
$ws.Range("K102").Value = 1723.3529
$ws.Range("L102").Value = 2257.4707
$ws.Range("M102").Value = -101.3529000000001
$ws.Range("N102").Value = -5501.4707
$ws.Range("H122").Value = 2419.8125
$ws.Range("I122").Value = 1246.2727
$ws.Range("K122").Value = 3738.8181
$ws.Range("M122").Value = -1288.8181
$ws.Range("H126").Value = 2932.12
$ws.Range("I126").Value = 2932.12
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8796.360000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -6326.360000000001
$ws.Range("H132").Value = 2613.3777
$ws.Range("I132").Value = 1677.0667
$ws.Range("J132").Value = 4486
$ws.Range("K132").Value = 5031.2001
$ws.Range("L132").Value = 13458
$ws.Range("M132").Value = -2501.2001
$ws.Range("N132").Value = -18518

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2151.9827
$ws.Range("I136").Value = 1363.561
$ws.Range("K136").Value = 4090.683
$ws.Range("M136").Value = -1540.683

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6078331
$ws.Range("I96").Value = 250949.5
$ws.Range("J96").Value = 11905712
$ws.Range("K96").Value = 250949.5
$ws.Range("L96").Value = 11905712
$ws.Range("M96").Value = -249576.5
$ws.Range("N96").Value = -11908458
$ws.Range("H126").Value = 249502.89
$ws.Range("I126").Value = 1498.591
$ws.Range("J126").Value = 509316.9
$ws.Range("K126").Value = 4495.772999999999
$ws.Range("L126").Value = 1527950.7
$ws.Range("M126").Value = -2025.772999999999
$ws.Range("N126").Value = -1532890.7
$ws.Range("H136").Value = 4797.022
$ws.Range("I136").Value = 5601.864
$ws.Range("J136").Value = 4059.25
$ws.Range("K136").Value = 16805.592
$ws.Range("L136").Value = 12177.75
$ws.Range("M136").Value = -14255.592
$ws.Range("N136").Value = -17277.75
